# Applies the LOQ4268.xlsx content update:
#  - Removes the old "Docentes responsáveis:" data row (old row 13, which held only the
#    professor name in columns B/C) so all subsequent rows shift up by one.
#  - Updates several cells to their new values, reflecting the refreshed course info.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row that used to hold only the professor's name in B/C (old row 13).
# This shifts every row below it up by one, matching the new layout (A1:C23).
$ws.Rows.Item(13).Delete()

# Row 10 (Objetivos / Objectives block): the descriptive objective text is replaced by the
# professor's identification text (string reused elsewhere in the sheet).
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

# Row 13 (Programa resumido:) now shows "Semestral" instead of the long summary text.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) now shows the activation date instead of the long program text.
# Copy the cell from the existing "Ativação:" row (B8/C8) so the date-like text
# "01/01/2021" is carried over as plain text (matching style/number format) instead of
# being auto-converted into a date serial number by Value assignment.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Row 18 (Método:) now shows the professor identification text.
$ws.Range("B18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C18").Value = "11079086 - Herlandí de Souza Andrade"

# Row 19 (Critério:) now shows the lecture method description.
$ws.Range("B19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."
$ws.Range("C19").Value = "Aulas expositivas teóricas, aulas práticas, aulas de exercícios."

# Row 20 (Norma de recuperação:) now shows the arithmetic-mean grading criterion text.
$ws.Range("B20").Value = "Média Aritmética das atividades avaliativas realizadas."
$ws.Range("C20").Value = "Média Aritmética das atividades avaliativas realizadas."

# Row 21 (Bibliografia:) now shows the recovery-grade norm text.
$ws.Range("B21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
$ws.Range("C21").Value = "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação."
